# Remove the "Image Directory" (H) and "File Name" (I) columns from the
# BookChapter sheet. The former column J ("URL") shifts left to become the
# new column H. Deleting the columns also drops the now-unreferenced shared
# strings (the "Image Directory"/"File Name" headers and the per-chapter
# image directory/file name values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("H:I").Delete()

# Reflect the last active cell/selection after the edit.
[void]$ws.Range("G3").Select()
